$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# --- Header row (row 1) updates ---
# N1: "Event" -> "Event.1"
$ws.Range("N1").Value = "Event.1"
# O1: "Correction " -> "Correction" (trailing space removed)
$ws.Range("O1").Value = "Correction"
# P1: new header "Serviced by " (trailing space kept)
$ws.Range("P1").Value = "Serviced by "
# copy the bold/bordered header formatting from O1 onto the new P1 header
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)

# --- Data rows 2-12 ---
# Column O was blank placeholder cells; they become literal "nan" text
# Column P is a brand-new blank placeholder column (mirrors the existing
# blank-cell convention used across the sheet, e.g. column O before this edit)
For ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = "nan"

    $ws.Cells.Item($r, 16).Value = ""
    $ws.Cells.Item($r, 16).Style = "Normal"
}

Write-Host "Card18 'Serviced by ' column added"
